$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) changes ---
$ws.Range("J1").Value2 = "official_notes"
$ws.Range("K1").Value2 = "researcher_notes"

# --- Data row (row 2) changes ---
$ws.Range("E2").Value2 = "Cumulative number of years of potential life lost from deaths among county residents under age 75 (summed over three years)."
$ws.Range("F2").Value2 = "Aggregate population under age 75 (over the three years) using bridged--race-postcensal estimates released by NCHS."
$ws.Range("J2").Value2 = "Data on deaths and births were provided by NCHS and drawn from the NVSS. These data are submitted to the NVSS by the vital registration systems operated in the jurisdictions legally responsible for registering vital events. Usually this variable is calculated by the NCHS, but more recently it has been calculated by the University of Wisconsin Population Health Institute directly using  the Mortality-All County micro-data."
$ws.Range("K2").Value2 = "Starting in 2020, each observation gets marked as unreliable or suppressed. A value of 0 means no flag is set, a value of 1 means the value is unreliable (based on death counts of 20 or less), and a value of 2 means the value was suppressed (9 or fewer deaths). As far as I can tell, earlier unreliable values are not reported as such."

# (New cells inherit the row's wrap-text / top-vertical alignment format
# automatically, matching the rest of the data row.)

# --- Row height (text now wraps across more lines) ---
$ws.Rows.Item(2).RowHeight = 114

# --- Column widths ---
# The underlying engine adds a constant +5/6 character offset when it
# serializes a column width back to OOXML, so we pre-compensate for it
# to land as close as possible to the target widths.
$offset = 5.0/6.0

$ws.Columns.Item(1).ColumnWidth = 15.8515625 - $offset
$ws.Columns.Item(5).ColumnWidth = 28.57421875 - $offset
$ws.Columns.Item(6).ColumnWidth = 29.57421875 - $offset
$ws.Columns.Item(8).ColumnWidth = 8.421875 - $offset
$ws.Columns.Item(9).ColumnWidth = 4.57421875 - $offset
$ws.Columns.Item(10).ColumnWidth = 50.00390625 - $offset
$ws.Columns.Item(11).ColumnWidth = 36.8515625 - $offset

Write-Output "edit applied"
